$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell for column C (row 1) ---
$ws.Range("C1").Value = "Avner 3.07"

# --- Row 3: Avner 3.07 results (PASS / green) ---
$ws.Range("A3").Value = "testSignIn, avnerg@perfectomobile.com, a1001a, Hi, Avner!"
$ws.Range("C3").Value = "PASS"
$ws.Range("C3").Interior.ColorIndex = 10
$ws.Range("D3").Value = "PASS"
$ws.Range("D3").Interior.ColorIndex = 10

# --- New header cell for column D (row 1) ---
$ws.Range("D1").Value = "Avner 3.08"

# --- Row 4: Avner 3.08 results (FAIL / red) ---
$ws.Range("A4").Value = "testSignIn, dudu@gulu.com, , Hi, Mister!"
$ws.Range("D4").Value = "FAIL"
$ws.Range("D4").Interior.ColorIndex = 3
$ws.Hyperlinks.Add($ws.Range("D4"), "C:/Users/rajp/git/Beton/Beton/test-output/screenshots/2015-08-10-16-18-37-EDT.png")

# --- Auto size the columns (as done in @AfterClass) ---
$ws.Columns.Item(1).ColumnWidth = 53.8
$ws.Columns.Item(2).ColumnWidth = 7.1
$ws.Columns.Item(3).ColumnWidth = 9.8
$ws.Columns.Item(4).ColumnWidth = 9.8

Write-Host "Edit applied"
